$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet: fix casing of ExerciceProfessionnel references ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("A6").Value = "DroitExerciceComplementaire.ExerciceProfessionnel"
$wsElem.Range("B6").Value = "DroitExerciceComplementaire.ExerciceProfessionnel"
$wsElem.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"
